$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 119.0815153333333
$ws.Cells.Item(2, 8).Value = 357.244546
$ws.Cells.Item(2, 9).Value = 0.431812569872284
$ws.Cells.Item(2, 10).Value = 0.4318125698722839
$ws.Cells.Item(2, 13).Value = 0.1419263333333333
$ws.Cells.Item(2, 14).Value = 0.425779
$ws.Cells.Item(2, 15).Value = 0.002583058778296354
$ws.Cells.Item(2, 16).Value = 0.002583058778296354
$ws.Cells.Item(2, 17).Value = 16.90080283903712
$ws.Cells.Item(2, 18).Value = 152.107225551334
$ws.Cells.Item(2, 19).Value = 0.001115397249187311
$ws.Cells.Item(2, 20).Value = 0.001115397249187311
$ws.Cells.Item(3, 7).Value = 119.0815153333333
$ws.Cells.Item(3, 8).Value = 357.244546
$ws.Cells.Item(3, 9).Value = 0.431812569872284
$ws.Cells.Item(3, 10).Value = 0.4318125698722839
$ws.Cells.Item(3, 15).Value = 0.001399682868699959
$ws.Cells.Item(3, 16).Value = 0.001399682868699959
$ws.Cells.Item(3, 17).Value = 9.15804332438689
$ws.Cells.Item(3, 18).Value = 82.422389919482
$ws.Cells.Item(3, 19).Value = 0.0006044006565395399
$ws.Cells.Item(3, 20).Value = 0.0006044006565395398
$ws.Cells.Item(4, 7).Value = 119.0815153333333
$ws.Cells.Item(4, 8).Value = 357.244546
$ws.Cells.Item(4, 9).Value = 0.431812569872284
$ws.Cells.Item(4, 10).Value = 0.4318125698722839
$ws.Cells.Item(4, 13).Value = 2.613991
$ws.Cells.Item(4, 14).Value = 7.841973
$ws.Cells.Item(4, 15).Value = 0.04757462720522382
$ws.Cells.Item(4, 16).Value = 0.04757462720522382
$ws.Cells.Item(4, 17).Value = 311.2780093476953
$ws.Cells.Item(4, 18).Value = 2801.502084129258
$ws.Cells.Item(4, 19).Value = 0.02054332203420357
$ws.Cells.Item(4, 20).Value = 0.02054332203420357
$ws.Cells.Item(5, 7).Value = 119.0815153333333
$ws.Cells.Item(5, 8).Value = 357.244546
$ws.Cells.Item(5, 9).Value = 0.431812569872284
$ws.Cells.Item(5, 10).Value = 0.4318125698722839
$ws.Cells.Item(5, 13).Value = 52.11224233333333
$ws.Cells.Item(5, 14).Value = 156.336727
$ws.Cells.Item(5, 15).Value = 0.9484426311477799
$ws.Cells.Item(5, 16).Value = 0.9484426311477798
$ws.Cells.Item(5, 17).Value = 6205.604784471217
$ws.Cells.Item(5, 18).Value = 55850.44306024094
$ws.Cells.Item(5, 19).Value = 0.4095494499323536
$ws.Cells.Item(5, 20).Value = 0.4095494499323535
$ws.Cells.Item(6, 9).Value = 0.4460879372303943
$ws.Cells.Item(6, 10).Value = 0.4460879372303942
$ws.Cells.Item(6, 13).Value = 0.1419263333333333
$ws.Cells.Item(6, 14).Value = 0.425779
$ws.Cells.Item(6, 15).Value = 0.002583058778296354
$ws.Cells.Item(6, 16).Value = 0.002583058778296354
$ws.Cells.Item(6, 17).Value = 17.459529439436
$ws.Cells.Item(6, 18).Value = 157.135764954924
$ws.Cells.Item(6, 19).Value = 0.001152271362155083
$ws.Cells.Item(6, 20).Value = 0.001152271362155083
$ws.Cells.Item(7, 9).Value = 0.4460879372303943
$ws.Cells.Item(7, 10).Value = 0.4460879372303942
$ws.Cells.Item(7, 15).Value = 0.001399682868699959
$ws.Cells.Item(7, 16).Value = 0.001399682868699959
$ws.Cells.Item(7, 19).Value = 0.0006243816436750855
$ws.Cells.Item(7, 20).Value = 0.0006243816436750855
$ws.Cells.Item(8, 9).Value = 0.4460879372303943
$ws.Cells.Item(8, 10).Value = 0.4460879372303942
$ws.Cells.Item(8, 13).Value = 2.613991
$ws.Cells.Item(8, 14).Value = 7.841973
$ws.Cells.Item(8, 15).Value = 0.04757462720522382
$ws.Cells.Item(8, 16).Value = 0.04757462720522382
$ws.Cells.Item(8, 17).Value = 321.568603563732
$ws.Cells.Item(8, 18).Value = 2894.117432073588
$ws.Cells.Item(8, 19).Value = 0.02122246731448329
$ws.Cells.Item(8, 20).Value = 0.02122246731448329
$ws.Cells.Item(9, 9).Value = 0.4460879372303943
$ws.Cells.Item(9, 10).Value = 0.4460879372303942
$ws.Cells.Item(9, 13).Value = 52.11224233333333
$ws.Cells.Item(9, 14).Value = 156.336727
$ws.Cells.Item(9, 15).Value = 0.9484426311477799
$ws.Cells.Item(9, 16).Value = 0.9484426311477798
$ws.Cells.Item(9, 17).Value = 6410.756959647068
$ws.Cells.Item(9, 18).Value = 57696.81263682361
$ws.Cells.Item(9, 19).Value = 0.4230888169100809
$ws.Cells.Item(9, 20).Value = 0.4230888169100808
$ws.Cells.Item(10, 7).Value = 33.50679633333333
$ws.Cells.Item(10, 8).Value = 100.520389
$ws.Cells.Item(10, 9).Value = 0.1215021138451521
$ws.Cells.Item(10, 10).Value = 0.121502113845152
$ws.Cells.Item(10, 13).Value = 0.1419263333333333
$ws.Cells.Item(10, 14).Value = 0.425779
$ws.Cells.Item(10, 15).Value = 0.002583058778296354
$ws.Cells.Item(10, 16).Value = 0.002583058778296354
$ws.Cells.Item(10, 17).Value = 4.755496745336778
$ws.Cells.Item(10, 18).Value = 42.799470708031
$ws.Cells.Item(10, 19).Value = 0.000313847101749283
$ws.Cells.Item(10, 20).Value = 0.0003138471017492829
$ws.Cells.Item(11, 7).Value = 33.50679633333333
$ws.Cells.Item(11, 8).Value = 100.520389
$ws.Cells.Item(11, 9).Value = 0.1215021138451521
$ws.Cells.Item(11, 10).Value = 0.121502113845152
$ws.Cells.Item(11, 15).Value = 0.001399682868699959
$ws.Cells.Item(11, 16).Value = 0.001399682868699959
$ws.Cells.Item(11, 17).Value = 2.576862509879222
$ws.Cells.Item(11, 18).Value = 23.191762588913
$ws.Cells.Item(11, 19).Value = 0.0001700644272598914
$ws.Cells.Item(11, 20).Value = 0.0001700644272598914
$ws.Cells.Item(12, 7).Value = 33.50679633333333
$ws.Cells.Item(12, 8).Value = 100.520389
$ws.Cells.Item(12, 9).Value = 0.1215021138451521
$ws.Cells.Item(12, 10).Value = 0.121502113845152
$ws.Cells.Item(12, 13).Value = 2.613991
$ws.Cells.Item(12, 14).Value = 7.841973
$ws.Cells.Item(12, 15).Value = 0.04757462720522382
$ws.Cells.Item(12, 16).Value = 0.04757462720522382
$ws.Cells.Item(12, 17).Value = 87.58646405416633
$ws.Cells.Item(12, 18).Value = 788.278176487497
$ws.Cells.Item(12, 19).Value = 0.005780417770829773
$ws.Cells.Item(12, 20).Value = 0.005780417770829772
$ws.Cells.Item(13, 7).Value = 33.50679633333333
$ws.Cells.Item(13, 8).Value = 100.520389
$ws.Cells.Item(13, 9).Value = 0.1215021138451521
$ws.Cells.Item(13, 10).Value = 0.121502113845152
$ws.Cells.Item(13, 13).Value = 52.11224233333333
$ws.Cells.Item(13, 14).Value = 156.336727
$ws.Cells.Item(13, 15).Value = 0.9484426311477799
$ws.Cells.Item(13, 16).Value = 0.9484426311477798
$ws.Cells.Item(13, 17).Value = 1746.114290336312
$ws.Cells.Item(13, 18).Value = 15715.0286130268
$ws.Cells.Item(13, 19).Value = 0.1152377845453131
$ws.Cells.Item(13, 20).Value = 0.1152377845453131
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.16474
$ws.Cells.Item(14, 8).Value = 0.49422
$ws.Cells.Item(14, 9).Value = 0.000597379052169715
$ws.Cells.Item(14, 10).Value = 0.000597379052169715
$ws.Cells.Item(14, 13).Value = 0.1419263333333333
$ws.Cells.Item(14, 14).Value = 0.425779
$ws.Cells.Item(14, 15).Value = 0.002583058778296354
$ws.Cells.Item(14, 16).Value = 0.002583058778296354
$ws.Cells.Item(14, 17).Value = 0.02338094415333334
$ws.Cells.Item(14, 18).Value = 0.21042849738
$ws.Cells.Item(14, 19).Value = [double]"1.543065204677338E-06"
$ws.Cells.Item(14, 20).Value = [double]"1.543065204677338E-06"
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.16474
$ws.Cells.Item(15, 8).Value = 0.49422
$ws.Cells.Item(15, 9).Value = 0.000597379052169715
$ws.Cells.Item(15, 10).Value = 0.000597379052169715
$ws.Cells.Item(15, 15).Value = 0.001399682868699959
$ws.Cells.Item(15, 16).Value = 0.001399682868699959
$ws.Cells.Item(15, 17).Value = 0.01266943952666667
$ws.Cells.Item(15, 18).Value = 0.11402495574
$ws.Cells.Item(15, 19).Value = [double]"8.361412254421692E-07"
$ws.Cells.Item(15, 20).Value = [double]"8.361412254421692E-07"
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.16474
$ws.Cells.Item(16, 8).Value = 0.49422
$ws.Cells.Item(16, 9).Value = 0.000597379052169715
$ws.Cells.Item(16, 10).Value = 0.000597379052169715
$ws.Cells.Item(16, 13).Value = 2.613991
$ws.Cells.Item(16, 14).Value = 7.841973
$ws.Cells.Item(16, 15).Value = 0.04757462720522382
$ws.Cells.Item(16, 16).Value = 0.04757462720522382
$ws.Cells.Item(16, 17).Value = 0.43062887734
$ws.Cells.Item(16, 18).Value = 3.87565989606
$ws.Cells.Item(16, 19).Value = [double]"2.842008570718414E-05"
$ws.Cells.Item(16, 20).Value = [double]"2.842008570718414E-05"
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 0.6666666666666666
$ws.Cells.Item(17, 7).Value = 0.16474
$ws.Cells.Item(17, 8).Value = 0.49422
$ws.Cells.Item(17, 9).Value = 0.000597379052169715
$ws.Cells.Item(17, 10).Value = 0.000597379052169715
$ws.Cells.Item(17, 13).Value = 52.11224233333333
$ws.Cells.Item(17, 14).Value = 156.336727
$ws.Cells.Item(17, 15).Value = 0.9484426311477799
$ws.Cells.Item(17, 16).Value = 0.9484426311477798
$ws.Cells.Item(17, 17).Value = 8.584970801993334
$ws.Cells.Item(17, 18).Value = 77.26473721794
$ws.Cells.Item(17, 19).Value = 0.0005665797600324115
$ws.Cells.Item(17, 20).Value = 0.0005665797600324114
